# Update countries & provincias Spain
# - Refresh case counts for a handful of countries (Dinamarca, Finlandia,
#   Uzbekistan, Kenia, Cabo Verde, Maldivas)
# - Update the "last updated" timestamp in A1
# - Re-sort the data table (A4:H216) by "Casos totales" (column B) descending,
#   since some of the updated totals change the country ranking order. The
#   newly-refreshed rows (Cabo Verde, Maldivas) are promoted ahead of any
#   other country they now tie with.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of country name -> updated column values (only columns that changed)
$updates = @{
    "Dinamarca"  = @{ D = 4312; E = 2839 }
    "Finlandia"  = @{ E = 2070; F = 67; G = 4; H = 98 }
    "Uzbekistan" = @{ D = 261; E = 1316 }
    "Kenia"      = @{ B = 281; C = 11; D = 69; E = 198 }
    "Cabo Verde" = @{ B = 67; C = 6; E = 65 }
    "Maldivas"   = @{ B = 60; C = 8; E = 44; F = 1 }
}

$colIndex = @{ A = 1; B = 2; C = 3; D = 4; E = 5; F = 6; G = 7; H = 8 }

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 4; $r -le $lastRow; $r++) {
    $name = $ws.Cells.Item($r, 1).Value2
    if ($updates.ContainsKey($name)) {
        $cols = $updates[$name]
        foreach ($col in $cols.Keys) {
            $c = $colIndex[$col]
            $ws.Cells.Item($r, $c).Value = $cols[$col]
        }
    }
}

# Update the "datos actualizados" timestamp
$ws.Range("A1").Value = "Datos actualizados a 20 de Abril de 2020 a las 14:52"

# Keep the table ranked by "Casos totales" (column B) descending
$dataRange = $ws.Range("A4:H" + $lastRow)
$keyRange = $ws.Range("B4:B" + $lastRow)
$dataRange.Sort($keyRange, 2)

# A stable sort leaves ties in their original relative order, but the data
# feed promotes a just-refreshed country ahead of whichever country it newly
# ties with. Fix up the two rows that just received a brand new value equal
# to their neighbour's.
function Swap-DataRows($ws, $r1, $r2) {
    for ($c = 1; $c -le 8; $c++) {
        $v1 = $ws.Cells.Item($r1, $c).Value2
        $v2 = $ws.Cells.Item($r2, $c).Value2
        $ws.Cells.Item($r1, $c).Value = $v2
        $ws.Cells.Item($r2, $c).Value = $v1
    }
}

for ($r = 5; $r -le $lastRow; $r++) {
    $name = $ws.Cells.Item($r, 1).Value2
    $rPrev = $r - 1
    if (($name -eq "Cabo Verde") -or ($name -eq "Maldivas")) {
        $bHere = $ws.Cells.Item($r, 2).Value2
        $bAbove = $ws.Cells.Item($rPrev, 2).Value2
        if ($bAbove -eq $bHere) {
            Swap-DataRows $ws $rPrev $r
        }
    }
}
